# Add USCDI Goals and Preferences part 1
# Two new "Observation" profiles are inserted into the "profiles" sheet,
# just above the existing "US Core Heart Rate Profile" row, matching the
# workbook's sort order (Type, then Name).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("profiles")

# Insert two blank rows at row 32 - everything that was row 32 and below
# (through the former row 60) shifts down to rows 34-62.
$ws.Rows.Item(32).Insert()
$ws.Rows.Item(32).Insert()

# New row 32: US Core Treatment Intervention Preference Profile
$ws.Range("A32").Value() = "http://hl7.org/fhir/us/core/StructureDefinition/us-core-treatment-intervention-preference"
$ws.Range("B32").Value() = "US Core Treatment Intervention Preference Profile"
$ws.Range("D32").Value() = "SHALL"
$ws.Range("E32").Value() = "Observation"

# New row 33: US Core Care Experience Preference Profile
$ws.Range("A33").Value() = "http://hl7.org/fhir/us/core/StructureDefinition/us-core-care-experience-preference"
$ws.Range("B33").Value() = "US Core Care Experience Preference Profile"
$ws.Range("D33").Value() = "SHALL"
$ws.Range("E33").Value() = "Observation"

# The two rows that previously carried the "just edited" highlight style
# (now at A44 - Average Blood Pressure, and A62 - the last row) revert to
# the plain/default style, since the highlight moved to the two brand new
# rows above.
$ws.Range("A44").Style = "Normal"
$ws.Range("A62").Style = "Normal"

# Update the view: scrolled near the top of the table with the last of the
# newly-sorted rows selected.
$ws.Range("A56").Select()
